# Edit script: insert 4 new weekly report rows at the top of the data
# (right before the existing row 22), pushing all the existing rows
# down by 4. This mirrors the commit "Fruta / hortaliza, semanal" which
# appends a new week's worth of price observations to the top of the
# historical log.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows starting at row 22; everything currently at row 22
# and below shifts down to row 26 and below.
$ws.Rows("22:25").Insert()

# Common boilerplate values shared by every data row in this sheet.
$mercadoId = 1
$mercado = "Agrícola del Norte S.A. de Arica"
$region = "Arica y Parinacota"
$codreg = 15
$tipo = "Fruta"
$productoId = 100104
$producto = "Frutos de pepita"
$categoriaId = 100104002
$categoria = "Manzana"

function Set-DataRow($Row, $Fecha, $Variedad, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $Unidad, $Origen, $PrecioKg, $KgUnidad) {
    $ws.Cells.Item($Row, 1).Value = $mercadoId
    $ws.Cells.Item($Row, 2).Value = $mercado
    $ws.Cells.Item($Row, 3).Value = $region
    $ws.Cells.Item($Row, 4).Value = $Fecha
    $ws.Cells.Item($Row, 5).Value = $codreg
    $ws.Cells.Item($Row, 6).Value = $tipo
    $ws.Cells.Item($Row, 7).Value = $productoId
    $ws.Cells.Item($Row, 8).Value = $producto
    $ws.Cells.Item($Row, 9).Value = $categoriaId
    $ws.Cells.Item($Row, 10).Value = $categoria
    $ws.Cells.Item($Row, 11).Value = $Variedad
    $ws.Cells.Item($Row, 12).Value = $Calidad
    $ws.Cells.Item($Row, 13).Value = $Volumen
    $ws.Cells.Item($Row, 14).Value = $PrecioMin
    $ws.Cells.Item($Row, 15).Value = $PrecioMax
    $ws.Cells.Item($Row, 16).Value = $PrecioProm
    $ws.Cells.Item($Row, 17).Value = $Unidad
    $ws.Cells.Item($Row, 18).Value = $Origen
    $ws.Cells.Item($Row, 19).Value = $PrecioKg
    $ws.Cells.Item($Row, 20).Value = $KgUnidad
}

Set-DataRow 22 44670 "Fuji royal"   "Segunda" 300 18000 19000 18500 "$/caja 20 kilos granel" "Región de O'Higgins" 925 20
Set-DataRow 23 44670 "Granny Smith" "Segunda" 300 18000 19000 18500 "$/caja 20 kilos granel" "Región de O'Higgins" 925 20
Set-DataRow 24 44670 "Royal Gala"   "Segunda" 270 18000 19000 18500 "$/caja 20 kilos granel" "Región de O'Higgins" 925 20
Set-DataRow 25 44670 "Scarlett"     "Segunda" 270 18000 19000 18500 "$/caja 20 kilos granel" "Región de O'Higgins" 925 20

# Match the date-cell number format used throughout column D.
$ws.Range("D22:D25").NumberFormat = $ws.Range("D21").NumberFormat
